# Update "paises" (countries) workbook:
#  - Insert "Sudan" into its correct alphabetical-ish slot (right after
#    Kirguistan / before Somalia) by re-pointing rows 97-100 to the
#    country names/stats that now occupy those positions (Sudan, Somalia,
#    Albania, Consejo Danes para los Refugiados), i.e. the former
#    "Sudan" entry (row 100) moves up to row 97 and everything that used
#    to sit between Kirguistan and Sudan shifts down by one row.
#  - Refresh a handful of per-country statistics (new confirmed cases /
#    deaths / etc.) for Estados Unidos, Alemania, Paises Bajos and Libia.
#  - Bump the "Datos actualizados" timestamp from 16:03 to 16:33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp header (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 16:33"

# --- plain statistic refreshes (country stays put) ----------------------
# Estados Unidos (row 4)
$ws.Range("B4").Value = 1239848
$ws.Range("C4").Value = 2215
$ws.Range("E4").Value = 966315
$ws.Range("G4").Value = 110
$ws.Range("H4").Value = 72381

# Alemania (row 9)
$ws.Range("B9").Value = 167372
$ws.Range("C9").Value = 365
$ws.Range("E9").Value = 22979

# Paises Bajos (row 19)
$ws.Range("F19").Value = 628

# Libia (row 170)
$ws.Range("B170").Value = 64
$ws.Range("C170").Value = 1
$ws.Range("E170").Value = 37

# --- Sudan inserted ahead of Somalia: rows 97-100 are rewritten in place
# row 97: Somalia -> Sudan (new data)
$ws.Range("A97").Value = "Sudan"
$ws.Range("B97").Value = 852
$ws.Range("C97").Value = 74
$ws.Range("D97").Value = 80
$ws.Range("E97").Value = 727
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 45

# row 98: Albania -> Somalia (former row-97 data)
$ws.Range("A98").Value = "Somalia"
$ws.Range("B98").Value = 835
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 75
$ws.Range("E98").Value = 722
$ws.Range("F98").Value = 2
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 38

# row 99: Consejo Danes para los Refugiados -> Albania (former row-98 data)
$ws.Range("A99").Value = "Albania"
$ws.Range("B99").Value = 832
$ws.Range("C99").Value = 12
$ws.Range("D99").Value = 595
$ws.Range("E99").Value = 206
$ws.Range("F99").Value = 7
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 31

# row 100: Sudan -> Consejo Danes para los Refugiados (former row-99 data)
$ws.Range("A100").Value = "Consejo Danes para los Refugiados"
$ws.Range("B100").Value = 797
$ws.Range("C100").Value = 92
$ws.Range("D100").Value = 92
$ws.Range("E100").Value = 670
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 35

# row 101 (Sri Lanka) is unaffected - left as-is.
